# Applies the "Add files via upload" edit:
#  - Sheet1 ("Sheet1" tab, xl/worksheets/sheet2.xml):
#      * rewrite C13's description text
#      * append 4 new test rows (14-17): upgrade / manda user / about / support
#      * widen column C to fit the new (longer) text
#      * leave the final selection on C17
#  - "orange login page" tab (xl/worksheets/sheet1.xml):
#      * scroll/zoom the view down to the bottom rows (55%) without changing
#        which sheet is active in the saved workbook

$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("orange login page")
$wsMain  = $wb.Worksheets.Item("Sheet1")

# --- Sheet1 ("Sheet1" tab): content edits -----------------------------------

# Existing row 13's description gets corrected/rewritten.
$wsMain.Range("C13").Value = "after click  to adime and then open new dashbood and then system user"

# New row 14 fills in the B/C cells for the already-present A14 (=201).
$wsMain.Range("B14").Value = "click on upgrade"
$wsMain.Range("C14").Value = "after click to upgrade then open to new url and new windowes open "

# New rows 15-17.
$wsMain.Range("A15").Value = 203
$wsMain.Range("B15").Value = "click on manda user"
$wsMain.Range("C15").Value = "after click to manda user and then open to new side box "

$wsMain.Range("A16").Value = 204
$wsMain.Range("B16").Value = "click on about "
$wsMain.Range("C16").Value = "after click to about then open to box then about informaction to companyname , version , avctive employess that all about "

$wsMain.Range("A17").Value = 205
$wsMain.Range("B17").Value = "click on support "
$wsMain.Range("C17").Value = "after click to support and then open to new page and then informaction to customer support "

# Column C needs to grow to fit the newly added (longer) descriptions.
$wsMain.Columns.Item(3).ColumnWidth = 112.5703125

# Last user action lands the selection on C17.
$wsMain.Activate()
$wsMain.Range("C17").Select()

# --- "orange login page" tab: view was scrolled down and zoomed out --------

$wsLogin.Activate()
$wsLogin.Range("A17").Select()
$excel.ActiveWindow.Zoom = 55

# Restore "Sheet1" as the active/visible tab (matches the saved workbook).
$wsMain.Activate()
